{"js": "const replacements = [\n  [\"290\u00f79=32, 2\", \"541\u00f78=67, 5\"],\n  [\"246\u00f78=30, 6\", \"784\u00f75=156, 4\"],\n  [\"713\u00f72=356, 1\", \"660\u00f79=73, 3\"],\n  [\"976\u00f73=325, 1\", \"664\u00f76=110, 4\"],\n  [\"217\u00f76=36, 1\", \"101\u00f77=14, 3\"],\n  [\"611\u00f75=122, 1\", \"781\u00f74=195, 1\"],\n  [\"682\u00f72=341, 0\", \"517\u00f75=103, 2\"],\n  [\"559\u00f74=139, 3\", \"352\u00f75=70, 2\"],\n  [\"766\u00f79=85, 1\", \"594\u00f77=84, 6\"],\n  [\"837\u00f76=139, 3\", \"587\u00f79=65, 2\"],\n  [\"274\u00f79=30, 4\", \"923\u00f74=230, 3\"],\n  [\"745\u00f78=93, 1\", \"211\u00f77=30, 1\"],\n  [\"225\u00f76=37, 3\", \"183\u00f74=45, 3\"],\n  [\"838\u00f74=209, 2\", \"923\u00f75=184, 3\"],\n  [\"178\u00f76=29, 4\", \"556\u00f76=92, 4\"],\n  [\"491\u00f73=163, 2\", \"208\u00f79=23, 1\"],\n  [\"816\u00f79=90, 6\", \"194\u00f79=21, 5\"],\n  [\"449\u00f72=224, 1\", \"740\u00f73=246, 2\"],\n  [\"651\u00f73=217, 0\", \"381\u00f77=54, 3\"],\n  [\"897\u00f77=128, 1\", \"903\u00f75=180, 3\"],\n  [\"496\u00f77=70, 6\", \"459\u00f72=229, 1\"],\n  [\"221\u00f76=36, 5\", \"269\u00f73=89, 2\"],\n  [\"565\u00f79=62, 7\", \"196\u00f72=98, 0\"],\n  [\"889\u00f77=127, 0\", \"159\u00f72=79, 1\"],\n  [\"998\u00f72=499, 0\", \"336\u00f74=84, 0\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"290\u00f79=32, 2\", \"541\u00f78=67, 5\"),\n  @(\"246\u00f78=30, 6\", \"784\u00f75=156, 4\"),\n  @(\"713\u00f72=356, 1\", \"660\u00f79=73, 3\"),\n  @(\"976\u00f73=325, 1\", \"664\u00f76=110, 4\"),\n  @(\"217\u00f76=36, 1\", \"101\u00f77=14, 3\"),\n  @(\"611\u00f75=122, 1\", \"781\u00f74=195, 1\"),\n  @(\"682\u00f72=341, 0\", \"517\u00f75=103, 2\"),\n  @(\"559\u00f74=139, 3\", \"352\u00f75=70, 2\"),\n  @(\"766\u00f79=85, 1\", \"594\u00f77=84, 6\"),\n  @(\"837\u00f76=139, 3\", \"587\u00f79=65, 2\"),\n  @(\"274\u00f79=30, 4\", \"923\u00f74=230, 3\"),\n  @(\"745\u00f78=93, 1\", \"211\u00f77=30, 1\"),\n  @(\"225\u00f76=37, 3\", \"183\u00f74=45, 3\"),\n  @(\"838\u00f74=209, 2\", \"923\u00f75=184, 3\"),\n  @(\"178\u00f76=29, 4\", \"556\u00f76=92, 4\"),\n  @(\"491\u00f73=163, 2\", \"208\u00f79=23, 1\"),\n  @(\"816\u00f79=90, 6\", \"194\u00f79=21, 5\"),\n  @(\"449\u00f72=224, 1\", \"740\u00f73=246, 2\"),\n  @(\"651\u00f73=217, 0\", \"381\u00f77=54, 3\"),\n  @(\"897\u00f77=128, 1\", \"903\u00f75=180, 3\"),\n  @(\"496\u00f77=70, 6\", \"459\u00f72=229, 1\"),\n  @(\"221\u00f76=36, 5\", \"269\u00f73=89, 2\"),\n  @(\"565\u00f79=62, 7\", \"196\u00f72=98, 0\"),\n  @(\"889\u00f77=127, 0\", \"159\u00f72=79, 1\"),\n  @(\"998\u00f72=499, 0\", \"336\u00f74=84, 0\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n\n  $found = $find.Execute(\n    $find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $find.Replacement.Text,\n    2\n  )\n\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}\n\nWrite-Output \"done\"\n"}
